$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 9 (remove extra expense entries), shifting cells up
$ws.Range("A4:C9").EntireRow.Delete()

# Update row 2: rent, 50000, 45782.125185185185
$ws.Range("A2").Value = "rent"
$ws.Range("B2").Value = 50000
$ws.Range("C2").Value = 45782.125185185185

# Update row 3: fuel, 20000, 45782.125185185185
$ws.Range("A3").Value = "fuel"
$ws.Range("B3").Value = 20000
$ws.Range("C3").Value = 45782.125185185185
